$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") '42.625.31'
Set-TextValue $ws.Range("E2") '  -0.44%  '

Set-TextValue $ws.Range("D3") '2.295.50'
Set-TextValue $ws.Range("E3") '  -0.19%  '

Set-TextValue $ws.Range("E4") '  -0.05%  '

Set-TextValue $ws.Range("D5") '301.19'
Set-TextValue $ws.Range("E5") '  -1.58%  '

Set-TextValue $ws.Range("D6") '95.56'
Set-TextValue $ws.Range("E6") '  -1.23%  '

Set-TextValue $ws.Range("E7") '  -0.45%  '

Set-TextValue $ws.Range("E9") '  -1.87%  '

Set-TextValue $ws.Range("D10") '34.54'
Set-TextValue $ws.Range("E10") '  -2.99%  '

Set-TextValue $ws.Range("D11") '19.19'
Set-TextValue $ws.Range("E11") '  +4.70%  '

Set-TextValue $ws.Range("E12") '  -0.99%  '

Set-TextValue $ws.Range("E13") '  -0.38%  '

Set-TextValue $ws.Range("D14") '6.75'
Set-TextValue $ws.Range("E14") '  +0.05%  '

Set-TextValue $ws.Range("D15") '2.649.22'
Set-TextValue $ws.Range("E15") '  -0.50%  '

Set-TextValue $ws.Range("D16") '2.308.24'
Set-TextValue $ws.Range("E16") '  +0.37%  '

Set-TextValue $ws.Range("D17") '0.782'
Set-TextValue $ws.Range("E17") '  +0.11%  '

Set-TextValue $ws.Range("D18") '42.543.13'
Set-TextValue $ws.Range("E18") '  -0.52%  '

Set-TextValue $ws.Range("D19") '12.28'
Set-TextValue $ws.Range("E19") '  -5.40%  '

Set-TextValue $ws.Range("E20") '  -1.08%  '

Set-TextValue $ws.Range("E21") '  -0.69%  '

Set-TextValue $ws.Range("D22") '67.80'
Set-TextValue $ws.Range("E22") '  +0.47%  '

Set-TextValue $ws.Range("D23") '2.27'
Set-TextValue $ws.Range("E23") '  +5.88%  '

Set-TextValue $ws.Range("D24") '234.94'
Set-TextValue $ws.Range("E24") '  -0.50%  '

Set-TextValue $ws.Range("E25") '  +0.16%  '

Set-TextValue $ws.Range("D26") '2.41'
Set-TextValue $ws.Range("E26") '  -2.50%  '

Set-TextValue $ws.Range("D27") '24.48'
Set-TextValue $ws.Range("E27") '  -3.56%  '

Set-TextValue $ws.Range("E28") '  +14.66%  '

Set-TextValue $ws.Range("D29") '164.56'
Set-TextValue $ws.Range("E29") '  -1.08%  '

Set-TextValue $ws.Range("D30") '9.05'
Set-TextValue $ws.Range("E30") '  -0.13%  '

Set-TextValue $ws.Range("D31") '32.13'
Set-TextValue $ws.Range("E31") '  -2.92%  '

Set-TextValue $ws.Range("E32") '  -0.02%  '

Set-TextValue $ws.Range("D33") '4.97'
Set-TextValue $ws.Range("E33") '  -0.41%  '

Set-TextValue $ws.Range("D34") '17.49'
Set-TextValue $ws.Range("E34") '  -1.32%  '

Set-TextValue $ws.Range("D35") '4.45'
Set-TextValue $ws.Range("E35") '  -6.75%  '

Set-TextValue $ws.Range("D36") '0.0703'
Set-TextValue $ws.Range("E36") '  +1.32%  '

Set-TextValue $ws.Range("E37") '  -3.18%  '

Set-TextValue $ws.Range("D38") '0.0999'
Set-TextValue $ws.Range("E38") '  -1.11%  '

Set-TextValue $ws.Range("E39") '  -0.53%  '

$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D40") '2.70'
Set-TextValue $ws.Range("E40") '  -0.37%  '

$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D41") '0.108'
Set-TextValue $ws.Range("E41") '  -1.14%  '

Set-TextValue $ws.Range("D42") '20.42'
Set-TextValue $ws.Range("E42") '  +12.00%  '

Set-TextValue $ws.Range("D43") '1.963.11'
Set-TextValue $ws.Range("E43") '  -1.99%  '

Set-TextValue $ws.Range("D44") '10.47'
Set-TextValue $ws.Range("E44") '  +4.81%  '

Set-TextValue $ws.Range("E45") '  -0.51%  '

Set-TextValue $ws.Range("E46") '  -2.33%  '

Set-TextValue $ws.Range("D47") '2.76'
Set-TextValue $ws.Range("E47") '  -0.43%  '

Set-TextValue $ws.Range("E48") '  -0.70%  '

Set-TextValue $ws.Range("D49") '2.521.66'
Set-TextValue $ws.Range("E49") '  -0.27%  '

Set-TextValue $ws.Range("D50") '53.12'
Set-TextValue $ws.Range("E50") '  -1.02%  '

Set-TextValue $ws.Range("D51") '71.29'
Set-TextValue $ws.Range("E51") '  -0.44%  '
